# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to match the refreshed data pull (gh-pages output regenerated
# at 456a3b4). Row numbers differ slightly between the two sheets
# because "全部类型" interleaves an extra row, so each sheet gets its
# own row => new-value map.

$wb = $excel.ActiveWorkbook

# row => new F-column value, for sheet "展览"
$updatesZhanLan = @{
    2  = 188
    4  = 148
    5  = 1306
    6  = 18216
    7  = 370
    10 = 6865
    11 = 690
    13 = 15
    15 = 65
    18 = 1303
    19 = 241
    21 = 657
    26 = 989
    28 = 5166
    30 = 39
    33 = 12097
    34 = 1284
    37 = 285
    38 = 3923
    40 = 91
}

# row => new F-column value, for sheet "全部类型"
$updatesQuanBu = @{
    2  = 188
    4  = 148
    5  = 1306
    6  = 18216
    7  = 370
    10 = 6865
    11 = 690
    13 = 15
    15 = 65
    18 = 1303
    19 = 241
    21 = 657
    26 = 989
    28 = 5166
    32 = 39
    35 = 12097
    36 = 1284
    39 = 285
    40 = 3923
    42 = 91
}

$wsZhanLan = $wb.Worksheets.Item("展览")
foreach ($row in $updatesZhanLan.Keys) {
    $wsZhanLan.Cells.Item($row, 6).Value = $updatesZhanLan[$row]
}

$wsQuanBu = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesQuanBu.Keys) {
    $wsQuanBu.Cells.Item($row, 6).Value = $updatesQuanBu[$row]
}
